# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1166
$ws1.Range("F3").Value = 51
$ws1.Range("F4").Value = 1451
$ws1.Range("F5").Value = 335
$ws1.Range("F6").Value = 1051
$ws1.Range("F7").Value = 10839
$ws1.Range("F8").Value = 28
$ws1.Range("F10").Value = 302
$ws1.Range("F11").Value = 1053
$ws1.Range("F12").Value = 731
$ws1.Range("F13").Value = 12150
$ws1.Range("F14").Value = 12613

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1166
$ws4.Range("F4").Value = 51
$ws4.Range("F5").Value = 1451
$ws4.Range("F6").Value = 335
$ws4.Range("F7").Value = 1051
$ws4.Range("F8").Value = 10839
$ws4.Range("F9").Value = 28
$ws4.Range("F11").Value = 302
$ws4.Range("F12").Value = 1053
$ws4.Range("F13").Value = 731
$ws4.Range("F14").Value = 12150
$ws4.Range("F15").Value = 12613
